$d = $word.ActiveDocument

# Step 1: Remove the "post - protected (user)" paragraph entirely (it is being
# dropped as part of generalizing the rating component).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*protected (user)*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# Step 2: Rewrite the description of rate_content so that it documents the
# generalized 1-10 rating scale instead of the old up/down rating scheme.
# Leave the leading "rate_content" word (and its spell-check markers) alone,
# and only replace the remainder of the sentence.
$oldText = " (/content/:contentUrl/submission/:submissionId/rating/:ratingVal): rates the piece of content with ratingVal (which must be " + [char]0x201C + "up" + [char]0x201D + " or " + [char]0x201C + "down" + [char]0x201D + ") "
$newText = " (protected) (/content/:contentUrl/submission/:submissionId/rating/:ratingVal): rating values are as follows: 1 => half star, 2 => 1 star, 3 => 1 " + [char]0x00BD + " star, " + [char]0x2026 + " and so on. The maximum rating is 10. "

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
